$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (shifts existing rows 4.. down by one)
$ws.Rows("4:4").Insert()

# Fill in the new benchmark entry: Ryzen 5 1600 AF
$ws.Range("A4").Value = "AMD"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = "Ryzen 5 1600 AF"
$ws.Range("D4").Value = 65
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 3.2
$ws.Range("H4").Value = 3.46
$ws.Range("I4").Value = "x86-64"
$ws.Range("J4").Value = "Performance"
$ws.Range("K4").Value = 32
$ws.Range("L4").Value = 2
$ws.Range("M4").Value = "DDR4"
$ws.Range("N4").Value = 3066
$ws.Range("O4").Value = 0.5
$ws.Range("P4").Value = 1.03
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 4.05

# Move selection to match post-edit cursor position
$ws.Range("R28").Select()
